$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set / update cell values ---
$ws.Range("B1").Value = "Mon Sep 19 2016 16:12:26 GMT+0700 (SE Asia Standard Time)"
$ws.Range("B2").Value = "3 sạch 17_09_2016"
$ws.Range("C2").Value = "HQ Food."
$ws.Range("C4").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C16").Value = 5
$ws.Range("C23").Value = 4
$ws.Range("C38").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("C42").Value = 5
$ws.Range("C44").Value = 1
$ws.Range("C45").Value = 3
$ws.Range("C62").Value = 4

# --- Clear obsolete cells ---
$clearAddrs = @("D2","B4","D5","B6","C6","B7","B8","B10","C12","B13","B16","D19","E19","B21","B22","B23","D23","C24","D24","D25","E25","F25","G25","H25","I25","B26","C26","B27","C27","B29","C29","B31","B33","B35","B36","B37","B38","B39","B40","B41","C41","B42","B43","B44","B45","B46","B47","C47","B48","B51","C51","B54","B55","B57","B58","B59","B60","B61","C61","B62","B63","B64")
foreach ($addr in $clearAddrs) {
    $ws.Range($addr).ClearContents()
}